# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 90 of Sheet1, pushing the
# existing rows 90-113 down to 91-114.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 90..113 down by one to make room for the new record.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new weekly data point.
$ws.Range("A90").Value = 3
$ws.Range("B90").Value = "Femacal de La Calera"
$ws.Range("C90").Value = "Coquimbo"
$ws.Range("D90").Value = 44551
$ws.Range("E90").Value = 5
$ws.Range("F90").Value = 100112030
$ws.Range("G90").Value = "Poroto granado"
$ws.Range("H90").Value = "Sin especificar"
$ws.Range("I90").Value = "Primera"
$ws.Range("J90").Value = 45
$ws.Range("K90").Value = 43000
$ws.Range("L90").Value = 45000
$ws.Range("M90").Value = 43889
$ws.Range("N90").Value = "$/malla 25 kilos"
$ws.Range("O90").Value = "Provincia de Limarí"
$ws.Range("P90").Value = 1756
$ws.Range("Q90").Value = 25
$ws.Range("R90").Value = "Hortaliza"
